# Implement cascade plot coordinates
# Adds an optional 'Plot Coordinates' column (C) to the "Compartments" sheet,
# giving finer control in plotting schematics.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compartments")

# Header for the new column
$ws.Range("C1").Value = "Plot Coordinates"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4108   # xlCenter

# Plot coordinates for each compartment row (2-31), centered alignment
$coords = @(
  "(-3,5)",
  "(3,5)",
  "(0,4)",
  "(-2,3)",
  "(-6,3)",
  "(2,3)",
  "(6,3)",
  "(-6,2)",
  "(-8.5,1)",
  "(-7.5,0)",
  "(-6.5,1)",
  "(-5.5,0)",
  "(-4.5,1)",
  "(-3.5,0)",
  "(0,2)",
  "(-2.5,1)",
  "(-1.5,0)",
  "(-0.5,1)",
  "(0.5,0)",
  "(1.5,1)",
  "(2.5,0)",
  "(6,2)",
  "(3.5,1)",
  "(4.5,0)",
  "(5.5,1)",
  "(6.5,0)",
  "(7.5,1)",
  "(8.5,0)",
  "(0,-1)",
  "(4,-1)"
)

for ($i = 0; $i -lt $coords.Length; $i++) {
  $row = $i + 2
  $cell = $ws.Cells.Item($row, 3)
  $cell.Value = $coords[$i]
  $cell.HorizontalAlignment = -4108   # xlCenter
}

# Size the new column to fit its contents
$ws.Columns.Item(3).ColumnWidth = 14.1666666666667

# The Compartments sheet becomes the active tab / selected sheet, with D30 selected
$ws.Activate() | Out-Null
$ws.Range("D30").Select() | Out-Null
